# Adding the changes we made on may 9th
#
# This inserts 3 new rows of data right after the header row (the
# existing rows 2-21 shift down to become rows 5-24), and appends
# 7 new rows of data after the old data block (new rows 25-31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 3 new rows right after the header row (before old row 2) ---
$ws.Rows("2:4").Insert()
# The insert above picks up the header row's formatting by default;
# the source data rows are unstyled, so strip that back off.
$ws.Rows("2:4").ClearFormats()

$newTopRows = @(
    @(-0.0041233403608202, 0.0484110713005065, -0.0001527163112768),
    @(-0.0736092627048492, -0.0381790772080421, 0.0797179117798805),
    @(-1.18019163608551, -4.37715482711792, 0.3266601860523224)
)

for ($i = 0; $i -lt $newTopRows.Count; $i++) {
    $r = 2 + $i
    $row = $newTopRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# --- Append 7 new rows of data after the (now shifted) old data block ---
# Old data originally occupied rows 2-21 (20 rows); after the insert above
# it now occupies rows 5-24. New rows go at 25-31.
$newBottomRows = @(
    @(-0.5545129179954529, -0.7066183686256409, -0.1945605874061584),
    @(-0.0233655963093042, -0.0335975885391235, -0.5940664410591125),
    @(0.1398881375789642, 0.0471893399953842, 0.531147301197052),
    @(-0.0551305897533893, 0.0639881342649459, 0.093156948685646),
    @(0.1805106848478317, 0.0415388382971286, 0.1635591685771942),
    @(-0.1348485052585601, 0.1539380401372909, 0.1916589736938476),
    @(-0.0897971913218498, 0.1873829066753387, -0.0282525178045034)
)

$startRow = 25
for ($i = 0; $i -lt $newBottomRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newBottomRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
